$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held two rows of data:
#   row 1: 2023-07-18 / 3525
#   row 2: 2023-07-19 / 4515
# The first row is being dropped entirely, and the second (surviving) row
# moves down to row 4, leaving rows 1-3 empty.

# Drop the first row (2023-07-18 / 3525); this shifts the remaining data
# (2023-07-19 / 4515) up into row 1.
$ws.Rows("1").Delete()

# Relocate that surviving record down onto row 4 (cut/paste keeps the date
# as literal text instead of re-parsing it into a date serial number).
$ws.Range("A1:B1").Cut($ws.Range("A4"))
